$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 17789
$ws.Range("J7").Value = 17862.5
$ws.Range("L7").Value = 17862.5
$ws.Range("N7").Value = -18086.5

$ws.Range("H12").Value = 161
$ws.Range("I12").Value = 161
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 161
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 9
$ws.Range("N12").ClearContents()

$ws.Range("H14").Value = 17789
$ws.Range("J14").Value = 17862.5
$ws.Range("L14").Value = 17862.5
$ws.Range("N14").Value = -18244.5

$ws.Range("H33").Value = 1003338.5
$ws.Range("I33").Value = 1502246.6
$ws.Range("K33").Value = 1502246.6
$ws.Range("M33").Value = -1502017.6

$ws.Range("H43").Value = 5250.231
$ws.Range("I43").Value = 6525.1665
$ws.Range("J43").Value = 4157.4287
$ws.Range("K43").Value = 6525.1665
$ws.Range("L43").Value = 4157.4287
$ws.Range("M43").Value = -6456.1665
$ws.Range("N43").Value = -4295.4287

$ws.Range("H64").Value = 5509.3335
$ws.Range("I64").Value = 3520
$ws.Range("K64").Value = 3520
$ws.Range("M64").Value = -3272

$ws.Range("H67").Value = 5509.3335
$ws.Range("I67").Value = 3520
$ws.Range("K67").Value = 3520
$ws.Range("M67").Value = -2662

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H138").Value = 3485.0833
$ws.Range("J138").Value = 3449.7727
$ws.Range("L138").Value = 10349.3181
$ws.Range("N138").Value = -20629.3181

$ws.Range("H141").Value = 2752.4736
$ws.Range("I141").Value = 2738.7222
$ws.Range("K141").Value = 8216.1666
$ws.Range("M141").Value = -3036.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8419.143
$ws.Range("I32").Value = 7235.5
$ws.Range("K32").Value = 7235.5
$ws.Range("M32").Value = -6948.5

$ws.Range("H97").Value = 1063.0286
$ws.Range("I97").Value = 1037.4062
$ws.Range("K97").Value = 1037.4062
$ws.Range("M97").Value = -541.4061999999999

$ws.Range("H122").Value = 5326.6665
$ws.Range("I122").Value = 3392
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 10176
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -7726
$ws.Range("N122").Value = -49900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5721.2085
$ws.Range("I20").Value = 3184.7058
$ws.Range("J20").Value = 11881.286
$ws.Range("K20").Value = 3184.7058
$ws.Range("L20").Value = 11881.286
$ws.Range("M20").Value = -2937.7058
$ws.Range("N20").Value = -12375.286

$ws.Range("H94").Value = 3023.4167
$ws.Range("I94").Value = 2920.111
$ws.Range("K94").Value = 2920.111
$ws.Range("M94").Value = -2469.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4511.75
$ws.Range("J31").Value = 5524.75
$ws.Range("L31").Value = 5524.75
$ws.Range("N31").Value = -6114.75

$ws.Range("H34").Value = 4511.75
$ws.Range("J34").Value = 5524.75
$ws.Range("L34").Value = 5524.75
$ws.Range("N34").Value = -5928.75

$ws.Range("H107").Value = 1618.8572
$ws.Range("I107").Value = 1888.2222
$ws.Range("J107").Value = 1134
$ws.Range("K107").Value = 1888.2222
$ws.Range("L107").Value = 1134
$ws.Range("M107").Value = 31.77780000000007
$ws.Range("N107").Value = -4974

$ws.Range("H132").Value = 21318.695
$ws.Range("I132").Value = 14195.131
$ws.Range("K132").Value = 42585.393
$ws.Range("M132").Value = -40055.393

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 58.055557
$ws.Range("I26").Value = 106.6
$ws.Range("J26").Value = 39.384617
$ws.Range("K26").Value = 319.8
$ws.Range("L26").Value = 118.153851
$ws.Range("M26").Value = -31.79999999999995
$ws.Range("N26").Value = -694.153851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 39359.375
$ws.Range("J92").Value = 39164.5
$ws.Range("L92").Value = 39164.5
$ws.Range("N92").Value = -42908.5

$ws.Range("H97").Value = 1386.2858
$ws.Range("I97").Value = 1386.2858
$ws.Range("K97").Value = 1386.2858
$ws.Range("M97").Value = -890.2858000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9696.637000000001
$ws.Range("I16").Value = 759.3158
$ws.Range("K16").Value = 759.3158
$ws.Range("M16").Value = -589.3158

$ws.Range("H46").Value = 1655.2354
$ws.Range("I46").Value = 1040.3334
$ws.Range("J46").Value = 1990.6364
$ws.Range("K46").Value = 1040.3334
$ws.Range("L46").Value = 1990.6364
$ws.Range("M46").Value = -852.3334
$ws.Range("N46").Value = -2366.6364

$ws.Range("H55").Value = 971.6667
$ws.Range("I55").Value = 624.5
$ws.Range("J55").Value = 1666
$ws.Range("K55").Value = 624.5
$ws.Range("L55").Value = 1666
$ws.Range("M55").Value = -451.5
$ws.Range("N55").Value = -2012

$ws.Range("H68").Value = 4064.2354
$ws.Range("I68").Value = 3006.125
$ws.Range("J68").Value = 5004.778
$ws.Range("K68").Value = 3006.125
$ws.Range("L68").Value = 5004.778
$ws.Range("M68").Value = -2257.125
$ws.Range("N68").Value = -6502.778

$ws.Range("H71").Value = 4064.2354
$ws.Range("I71").Value = 3006.125
$ws.Range("J71").Value = 5004.778
$ws.Range("K71").Value = 15030.625
$ws.Range("L71").Value = 25023.89
$ws.Range("M71").Value = -11286.625
$ws.Range("N71").Value = -32511.89

$ws.Range("H93").Value = 4259
$ws.Range("I93").Value = 4080.756
$ws.Range("J93").Value = 4868
$ws.Range("K93").Value = 4080.756
$ws.Range("L93").Value = 4868
$ws.Range("M93").Value = -2832.756
$ws.Range("N93").Value = -7364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 9497.861000000001
$ws.Range("I126").Value = 6242.3447
$ws.Range("J126").Value = 22985
$ws.Range("K126").Value = 18727.0341
$ws.Range("L126").Value = 68955
$ws.Range("M126").Value = -16257.0341
$ws.Range("N126").Value = -73895
